$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- xs-TAINI-laser1: add checks for non-numeric subject IDs --------------
# The per-subject "is this box in use" indicators (row 9) used to test
# whether the SUBJECT-ID cell held a positive number. Subject IDs are not
# always numeric, so the check now simply treats a blank SUBJECT-ID box as
# "unused" and anything else (numeric or not) as "used".
$ws.Range("C9").Formula = "=IF(SUBJECT1=`"`",0,1)"
$ws.Range("E9").Formula = "=IF(SUBJECT2=`"`",0,1)"
$ws.Range("G9").Formula = "=IF(SUBJECT3=`"`",0,1)"
$ws.Range("I9").Formula = "=IF(SUBJECT4=`"`",0,1)"

# LBOX1 (C10) loses its placeholder value of 1 - the box numbering for
# subject 1 (column B, rows 15-105) now starts counting from an empty
# (0-valued) LBOX1, shifting every cached box number down by one.
$ws.Range("C10").Value = ""

# --- cosmetic workbook/view state also touched by this save ---------------
$ws.StandardWidth = 9.19140625
$ws.Range("D13").Select()
